# Insert a new weekly record as row 15, pushing existing rows 15-28 down to 16-29.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(15).Insert()

$ws.Cells.Item(15, 1).Value = 9
$ws.Cells.Item(15, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(15, 3).Value = "Metropolitana"
$ws.Cells.Item(15, 4).Value = "9/8/2021"
$ws.Cells.Item(15, 5).Value = 13
$ws.Cells.Item(15, 6).Value = 100112029
$ws.Cells.Item(15, 7).Value = "Orégano"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 16
$ws.Cells.Item(15, 11).Value = 10000
$ws.Cells.Item(15, 12).Value = 10500
$ws.Cells.Item(15, 13).Value = 10250
$ws.Cells.Item(15, 14).Value = "`$/docena de atados"
$ws.Cells.Item(15, 15).Value = "Región Metropolitana"
$ws.Cells.Item(15, 16).Value = 3417
$ws.Cells.Item(15, 17).Value = 3
$ws.Cells.Item(15, 18).Value = "Hortaliza"
